$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "General"
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("B3").Value = 41096.520833333336

# ---------------------------------------------------------------------------
# Sheet "BESS"
# ---------------------------------------------------------------------------
$wsBess = $wb.Worksheets.Item("BESS")

# Bold the "kV" header
$wsBess.Range("D1").Font.Bold = $true

# New value in C2
$wsBess.Range("C2").Value = 3

# Updated numbers on row 2
$wsBess.Range("E2").Value = 10
$wsBess.Range("F2").Value = 15
$wsBess.Range("G2").Value = 30

# Row 3 data removed entirely (keeps formatting on D3/H3/I3)
$wsBess.Range("A3:I3").ClearContents()

[void]$wsBess.Range("C6").Select()

# ---------------------------------------------------------------------------
# Sheet "Generator"
# ---------------------------------------------------------------------------
$wsGenerator = $wb.Worksheets.Item("Generator")

$wsGenerator.Range("G2").Value = 1
$wsGenerator.Range("H2").Value = 1
$wsGenerator.Range("I2").Value = "dados_power_ajustado"

$wsGenerator.Columns.Item(9).ColumnWidth = 18.75

[void]$wsGenerator.Range("G4").Select()

# ---------------------------------------------------------------------------
# Sheet "Load"
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")

$wsLoad.Range("F1").Value = "Pmax"

$wsLoad.Range("B3").Value = 2
$wsLoad.Range("B4").Value = 3

$wsLoad.Range("K2").Value = 12
$wsLoad.Range("K3").Value = 12

$wsLoad.Range("D8").Font.Underline = $true
$wsLoad.Range("I9").Font.Underline = $true

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping (matches final workbook view state)
# ---------------------------------------------------------------------------
[void]$wsGeneral.Activate()
[void]$wsGeneral.Range("D6").Select()
